$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three obsolete worker/period detail rows (MARIA PATRICIA POLO PEREIRA
# rows and the older KEREN ESTHER CANOLES PEREIRA period row), keeping the last
# detail row (KEREN ESTHER CANOLES PEREIRA) which becomes the sole remaining entry.
$ws.Rows("16:18").Delete()

# The surviving detail row now reports the new "2508" period (part 1 of the new
# account statement) instead of the old "2506" period.
$ws.Range("E16").Value = "2508"

# Update the summary figures to reflect a single worker / single period now.
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
